$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2,4).Value = "29.296.58"
$ws.Cells.Item(2,5).Value = "  -0.37%  "

# Row 3: Ethereum
$ws.Cells.Item(3,4).Value = "1.844.29"
$ws.Cells.Item(3,5).Value = "  -0.35%  "

# Row 4: TetherUSD
$ws.Cells.Item(4,4).Value = "'1.003"
$ws.Cells.Item(4,5).Value = "  +0.35%  "

# Row 5: BNB
$ws.Cells.Item(5,4).Value = "'239.80"
$ws.Cells.Item(5,5).Value = "  -0.43%  "

# Row 6: XRP
$ws.Cells.Item(6,4).Value = "'0.6264"
$ws.Cells.Item(6,5).Value = "  -0.40%  "

# Row 7: USDC
$ws.Cells.Item(7,4).Value = "'0.9977"
$ws.Cells.Item(7,5).Value = "  -0.27%  "

# Row 8: Dogecoin
$ws.Cells.Item(8,4).Value = "'0.07593"
$ws.Cells.Item(8,5).Value = "  -1.04%  "

# Row 9: Cardano
$ws.Cells.Item(9,4).Value = "'0.2908"
$ws.Cells.Item(9,5).Value = "  -0.87%  "

# Row 10: Solana
$ws.Cells.Item(10,4).Value = "'24.54"
$ws.Cells.Item(10,5).Value = "  +0.00%  "

# Row 11: TRON
$ws.Cells.Item(11,4).Value = "'0.07729"
$ws.Cells.Item(11,5).Value = "  -0.24%  "

# Row 12: Polkadot
$ws.Cells.Item(12,4).Value = "'5.006"
$ws.Cells.Item(12,5).Value = "  -0.29%  "

# Row 13: Polygon
$ws.Cells.Item(13,4).Value = "'0.6779"
$ws.Cells.Item(13,5).Value = "  -0.48%  "

# Row 14: ShibaInu
$ws.Cells.Item(14,5).Value = "  -5.05%  "

# Row 15: Litecoin
$ws.Cells.Item(15,4).Value = "'83.00"
$ws.Cells.Item(15,5).Value = "  -0.74%  "

# Row 16: Uniswap
$ws.Cells.Item(16,4).Value = "'6.112"
$ws.Cells.Item(16,5).Value = "  -0.55%  "

# Row 17: WrappedBTC
$ws.Cells.Item(17,4).Value = "29.317.17"
$ws.Cells.Item(17,5).Value = "  -0.46%  "

# Row 18: BitcoinCash
$ws.Cells.Item(18,4).Value = "'228.72"
$ws.Cells.Item(18,5).Value = "  -0.22%  "

# Row 19: Avalanche
$ws.Cells.Item(19,4).Value = "'12.30"
$ws.Cells.Item(19,5).Value = "  -1.26%  "

# Row 20: Dai
$ws.Cells.Item(20,4).Value = "'0.9985"
$ws.Cells.Item(20,5).Value = "  -0.17%  "

# Row 21: Chainlink
$ws.Cells.Item(21,4).Value = "'7.458"
$ws.Cells.Item(21,5).Value = "  +0.27%  "

# Row 22: BinanceUSD
$ws.Cells.Item(22,4).Value = "'0.9988"
$ws.Cells.Item(22,5).Value = "  -0.19%  "

# Row 23: Monero
$ws.Cells.Item(23,4).Value = "'158.14"
$ws.Cells.Item(23,5).Value = "  +0.77%  "

# Row 24: Stellar
$ws.Cells.Item(24,4).Value = "'0.1386"
$ws.Cells.Item(24,5).Value = "  +0.10%  "

# Row 25: Cosmos
$ws.Cells.Item(25,4).Value = "'8.418"
$ws.Cells.Item(25,5).Value = "  +0.35%  "

# Row 26: EthereumClassic
$ws.Cells.Item(26,4).Value = "'17.64"
$ws.Cells.Item(26,5).Value = "  -0.24%  "

# Row 27: Toncoin
$ws.Cells.Item(27,4).Value = "'1.454"
$ws.Cells.Item(27,5).Value = "  +10.41%  "

# Row 28: PancakeSwap
$ws.Cells.Item(28,4).Value = "'1.471"
$ws.Cells.Item(28,5).Value = "  +0.32%  "

# Row 29: Hedera
$ws.Cells.Item(29,4).Value = "'0.05588"
$ws.Cells.Item(29,5).Value = "  -2.10%  "

# Row 30: Filecoin
$ws.Cells.Item(30,4).Value = "'4.097"
$ws.Cells.Item(30,5).Value = "  -0.68%  "

# Row 31: InternetComputer(DFINITY)
$ws.Cells.Item(31,4).Value = "'4.066"
$ws.Cells.Item(31,5).Value = "  +0.39%  "

# Row 32: LidoDAOToken
$ws.Cells.Item(32,4).Value = "'1.828"
$ws.Cells.Item(32,5).Value = "  -1.11%  "

# Row 33: ARBITRUM
$ws.Cells.Item(33,4).Value = "'1.153"
$ws.Cells.Item(33,5).Value = "  -0.91%  "

# Row 34: ImmutableX
$ws.Cells.Item(34,4).Value = "'0.6998"
$ws.Cells.Item(34,5).Value = "  -0.62%  "

# Row 35: HuobiToken
$ws.Cells.Item(35,4).Value = "'2.580"
$ws.Cells.Item(35,5).Value = "  -0.12%  "

# Row 36: Maker
$ws.Cells.Item(36,4).Value = "1.233.15"
$ws.Cells.Item(36,5).Value = "  +1.19%  "

# Row 37: VeChain
$ws.Cells.Item(37,4).Value = "'0.01803"
$ws.Cells.Item(37,5).Value = "  +0.72%  "

# Row 38: MXToken
$ws.Cells.Item(38,4).Value = "'2.728"
$ws.Cells.Item(38,5).Value = "  -1.88%  "

# Row 39: FraxShare
$ws.Cells.Item(39,4).Value = "'6.409"
$ws.Cells.Item(39,5).Value = "  -1.36%  "

# Row 40: TrustWalletToken
$ws.Cells.Item(40,4).Value = "'0.9039"
$ws.Cells.Item(40,5).Value = "  -0.54%  "

# Row 41: PaxDollar
$ws.Cells.Item(41,4).Value = "'0.9974"
$ws.Cells.Item(41,5).Value = "  -0.28%  "

# Row 42: Quant
$ws.Cells.Item(42,4).Value = "'101.30"
$ws.Cells.Item(42,5).Value = "  -0.46%  "

# Row 43: Aave
$ws.Cells.Item(43,4).Value = "'65.37"
$ws.Cells.Item(43,5).Value = "  -1.42%  "

# Row 44: Aptos
$ws.Cells.Item(44,4).Value = "'7.172"
$ws.Cells.Item(44,5).Value = "  +0.80%  "

# Row 45: TheSandbox
$ws.Cells.Item(45,4).Value = "'0.3985"
$ws.Cells.Item(45,5).Value = "  -0.73%  "

# Row 46: Algorand
$ws.Cells.Item(46,2).Value = "Algorand"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(46,4).Value = "'0.1152"
$ws.Cells.Item(46,5).Value = "  +1.98%  "

# Row 47: EnergySwap
$ws.Cells.Item(47,2).Value = "EnergySwap"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47,4).Value = "'8.979"
$ws.Cells.Item(47,5).Value = "  -0.01%  "

# Row 48: RenderToken
$ws.Cells.Item(48,4).Value = "'1.682"
$ws.Cells.Item(48,5).Value = "  +0.11%  "

# Row 49: BabyDogeCoin
$ws.Cells.Item(49,4).Value = "'0.00000000113"
$ws.Cells.Item(49,5).Value = "  -6.65%  "

# Row 50: Cronos
$ws.Cells.Item(50,4).Value = "'0.05695"
$ws.Cells.Item(50,5).Value = "  -0.27%  "

# Row 51: Mantle
$ws.Cells.Item(51,4).Value = "'0.4615"
$ws.Cells.Item(51,5).Value = "  -0.30%  "
